$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule" (sheet1): update Cost ($) and Unit Cost ($/ML) on row 2 ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 830.7307515000005
$wsSchedule.Range("F2").Value = 13.73562750496033

# --- Sheet "Detailed" (sheet2): update Price column (B) and some Type values (C) ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B19").Value = 32.18984

$wsDetailed.Range("B20").Value = 0.7

$wsDetailed.Range("B21").Value = -7.5187
$wsDetailed.Range("C21").Value = "historical"

$wsDetailed.Range("B22").Value = -8.42065
$wsDetailed.Range("C22").Value = "historical"

$wsDetailed.Range("B23").Value = -8.49024

$wsDetailed.Range("B24").Value = -6.83333

$wsDetailed.Range("B25").Value = -8.40512

$wsDetailed.Range("B26").Value = 0.01766

$wsDetailed.Range("B27").Value = 0.0263

$wsDetailed.Range("B28").Value = 0.02612

$wsDetailed.Range("B29").Value = -6.83333

$wsDetailed.Range("B30").Value = -5.58973

$wsDetailed.Range("B31").Value = -5.62882

$wsDetailed.Range("B32").Value = -6.28212

$wsDetailed.Range("B33").Value = -4.38967

$wsDetailed.Range("B34").Value = -1.07421

$wsDetailed.Range("B35").Value = -4.36085

$wsDetailed.Range("B36").Value = 0

$wsDetailed.Range("B37").Value = 33.14118

$wsDetailed.Range("B38").Value = 48.42444

$wsDetailed.Range("B39").Value = 56.98

$wsDetailed.Range("B40").Value = 57.06

$wsDetailed.Range("B41").Value = 58.86381

$wsDetailed.Range("B42").Value = 58.07672

$wsDetailed.Range("B44").Value = 57.09
